$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# Remove the "Caminho para pasta do cliente" column (G) entirely.
$ws.Range("G1:G2").EntireColumn.Delete()

# Update the visible selection to match the post-edit state.
$ws.Range("A3:XFD11").Select()
